$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# Row 1: "optimization_parameter" / "value" header row loses its trailing
# duplicated "value" cells (C1:F1) -- only A1/B1 remain.
$ws.Range("C1:F1").ClearContents()

# Row 8: the "Model" label is renamed to "production_function" (the
# "Sigmoid" value in B8 is unchanged).
$ws.Range("A8").Value = "production_function"
$ws.Range("B8").Value = "Sigmoid"

# A new row is inserted right after row 8 ("L_curve" = 0), pushing the
# previously-following rows down by one. Rather than using a native
# row-insert (which perturbs unrelated floating point literals elsewhere
# on the sheet), every affected row is rewritten explicitly below.

# New row 9: L_curve
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 0
$ws.Range("B9").Style = "Comma [0]"
$ws.Range("B9").NumberFormat = "0.00E+00"

# Row 10 (was row 9): estimate_params
$ws.Range("A10").Value = "estimate_params"
$ws.Range("B10").Value = 1

# Row 11 (was row 10): make_graphs
$ws.Range("A11").Value = "make_graphs"
$ws.Range("B11").Value = 0

# Row 12 (was row 11): fix_P
$ws.Range("A12").Value = "fix_P"
$ws.Range("B12").Value = 1

# Row 13 (was row 12): fix_b
$ws.Range("A13").Value = "fix_b"
$ws.Range("B13").Value = 1

# Row 14 (was row 13): expression_timepoints
$ws.Range("A14").Value = "expression_timepoints"
$ws.Range("B14").Value = 0.4
$ws.Range("C14").Value = 0.8
$ws.Range("D14").Value = 1.2
$ws.Range("E14").Value = 1.6
$ws.Range("B14").NumberFormat = "General"

# Row 15 (was row 14): Strain
$ws.Range("A15").Value = "Strain"
$ws.Range("B15").Value = "wt"
$ws.Range("C15").Value = "dcin5"
$ws.Range("B15").NumberFormat = "General"

# Row 16 (was row 15, the old "Deletion" row 16 is dropped entirely):
# Sheet
$ws.Range("A16").Value = "Sheet"
$ws.Range("B16").Value = 3
$ws.Range("C16").Value = 4

# Row 17 (simulation_timepoints) is untouched -- its row number and every
# B17:V17 numeric literal stay exactly as they were.

# The active tab moves from "network_weights" to "optimization_parameters",
# and the latter's selection becomes C1:H3.
$ws.Activate()
$ws.Range("C1:H3").Select()
